$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Delete the "Meta description: ..." paragraph near the top of the
#    document (it moves down near the end of the doc, see step 3).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$null = $metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Replace the final "Prompt: DALLE..." paragraph's text with the old
#    meta-description body text (keeps the existing italic run).
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Prompt: DALLE*slot game.",
    $false, $false, $true, $false, $false,
    $true, 1, $false,
    "Experience the Montreal circus-themed Cirque Du Soleil Kooza slot game for free. Read our review covering gameplay, graphics, theme, and pro and cons.",
    2)

# ------------------------------------------------------------------
# 3. Insert a new bold paragraph with the page title just before that
#    final paragraph.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$anchorPara = $d.Paragraphs.Item($n - 1)
$null = $anchorPara.Range.InsertParagraphAfter()
$insertTarget = $d.Paragraphs.Item($n).Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes" ?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cirque Du Soleil Kooza for Free: Review and Features</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $insertTarget.InsertXML($xml)
